$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.792708522335772
$ws.Range("C2").Value = 0.2544466942002543
$ws.Range("D2").Value = 0.08884395901115028
$ws.Range("E2").Value = 0.04573766239179378
$ws.Range("F2").Value = 2.122403679451466
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 1.546942915479946
$ws.Range("L2").Value = 0.2749935328722728

$ws.Range("B3").Value = 1.678671900351219
$ws.Range("C3").Value = 0.2219482845710559
$ws.Range("D3").Value = 0.08912226837352932
$ws.Range("E3").Value = 0.045893536406727
$ws.Range("F3").Value = 2.068054950713929
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 1.525602533049593
$ws.Range("L3").Value = 0.2647931738953702

$ws.Range("B4").Value = 1.60975129926544
$ws.Range("C4").Value = 0.2020316225883505
$ws.Range("D4").Value = 0.08932779141307989
$ws.Range("E4").Value = 0.04599787337165928
$ws.Range("F4").Value = 2.036019618066163
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 1.513369054355991
$ws.Range("L4").Value = 0.2587060589917769

$ws.Range("B5").Value = 1.581940629924816
$ws.Range("C5").Value = 0.1939241500998321
$ws.Range("D5").Value = 0.08942018263933704
$ws.Range("E5").Value = 0.04604256283596664
$ws.Range("F5").Value = 2.023297626148334
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 1.508600719692993
$ws.Range("L5").Value = 0.256269520896069

$ws.Range("B6").Value = 1.577339268476862
$ws.Range("C6").Value = 0.1925784175071215
$ws.Range("D6").Value = 0.08943604381196479
$ws.Range("E6").Value = 0.04605011466805475
$ws.Range("F6").Value = 2.021205153028546
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 1.507821990659153
$ws.Range("L6").Value = 0.2558675870447189

$ws.Range("B7").Value = 1.609375122740062
$ws.Range("C7").Value = 0.2019222479424911
$ws.Range("D7").Value = 0.08932900255081933
$ws.Range("E7").Value = 0.04599846727586954
$ws.Range("F7").Value = 2.035846701905996
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 1.51330387102584
$ws.Range("L7").Value = 0.2586730210914396

$ws.Range("B8").Value = 1.753159983072692
$ws.Range("C8").Value = 0.2432330273925913
$ws.Range("D8").Value = 0.08893268737354987
$ws.Range("E8").Value = 0.0457896174260628
$ws.Range("F8").Value = 2.103385306817657
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 1.539403194276758
$ws.Range("L8").Value = 0.2714397861484201

$ws.Range("B9").Value = 2.043907875067134
$ws.Range("C9").Value = 0.3245766972718798
$ws.Range("D9").Value = 0.08843361450765741
$ws.Range("E9").Value = 0.04544848323486361
$ws.Range("F9").Value = 2.246566849252986
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 1.59756731431267
$ws.Range("L9").Value = 0.2978837439293045

$ws.Range("B10").Value = 2.263003415741593
$ws.Range("C10").Value = 0.3845990521124918
$ws.Range("D10").Value = 0.08824100734488383
$ws.Range("E10").Value = 0.04523949753715772
$ws.Range("F10").Value = 2.358523981160914
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 1.644679326440013
$ws.Range("L10").Value = 0.3181903824846728

$ws.Range("B11").Value = 2.363896293203652
$ws.Range("C11").Value = 0.4119736977110051
$ws.Range("D11").Value = 0.08819216580236855
$ws.Range("E11").Value = 0.04515345536422366
$ws.Range("F11").Value = 2.410971984475196
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 1.667089367148805
$ws.Range("L11").Value = 0.3276235987049034

$ws.Range("B12").Value = 2.40227999684987
$ws.Range("C12").Value = 0.4223507785994798
$ws.Range("D12").Value = 0.08817932775094306
$ws.Range("E12").Value = 0.04512217060476686
$ws.Range("F12").Value = 2.431054634634108
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 1.675718259728114
$ws.Range("L12").Value = 0.3312241668819667

$ws.Range("B13").Value = 2.394005447673635
$ws.Range("C13").Value = 0.4201153909153277
$ws.Range("D13").Value = 0.08818183975051141
$ws.Range("E13").Value = 0.0451288506354004
$ws.Range("F13").Value = 2.426719560200979
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 1.67385349508416
$ws.Range("L13").Value = 0.3304474518009926

$ws.Range("B14").Value = 2.367050571730886
$ws.Range("C14").Value = 0.4128272033370308
$ws.Range("D14").Value = 0.08819099576256662
$ws.Range("E14").Value = 0.04515085554580356
$ws.Range("F14").Value = 2.412619732660175
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 1.667796400064063
$ws.Range("L14").Value = 0.3279192485101703

$ws.Range("B15").Value = 2.350563129600459
$ws.Range("C15").Value = 0.4083644212847162
$ws.Range("D15").Value = 0.08819734323812867
$ws.Range("E15").Value = 0.04516450315501519
$ws.Range("F15").Value = 2.404012159129309
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 1.664104896914964
$ws.Range("L15").Value = 0.3263743598546966

$ws.Range("B16").Value = 2.256434616480419
$ws.Range("C16").Value = 0.3828115352207533
$ws.Range("D16").Value = 0.08824498645901002
$ws.Range("E16").Value = 0.04524530219119782
$ws.Range("F16").Value = 2.355127221408964
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 1.643234637527385
$ws.Range("L16").Value = 0.3175778601506209

$ws.Range("B17").Value = 2.19900472838242
$ws.Range("C17").Value = 0.3671542086014483
$ws.Range("D17").Value = 0.08828420313603402
$ws.Range("E17").Value = 0.04529718133620086
$ws.Range("F17").Value = 2.325529151467663
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 1.630683444075146
$ws.Range("L17").Value = 0.3122317960715577

$ws.Range("B18").Value = 2.166087800250011
$ws.Range("C18").Value = 0.3581550844534718
$ws.Range("D18").Value = 0.08831040445686966
$ws.Range("E18").Value = 0.04532787048554798
$ws.Range("F18").Value = 2.308647812362182
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 1.623556315742505
$ws.Range("L18").Value = 0.309175277152832

$ws.Range("B19").Value = 2.154962437060078
$ws.Range("C19").Value = 0.3551092373846814
$ws.Range("D19").Value = 0.08831989933131013
$ws.Range("E19").Value = 0.04533840724490723
$ws.Range("F19").Value = 2.302956484099127
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 1.621158931153417
$ws.Range("L19").Value = 0.3081435433301607

$ws.Range("B20").Value = 2.205106301525063
$ws.Range("C20").Value = 0.368820273447227
$ws.Range("D20").Value = 0.08827965063227339
$ws.Range("E20").Value = 0.04529157078185797
$ws.Range("F20").Value = 2.328665124218901
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 1.632010005528258
$ws.Range("L20").Value = 0.3127989877992974

$ws.Range("B21").Value = 2.374963034897917
$ws.Range("C21").Value = 0.4149676184583768
$ws.Range("D21").Value = 0.08818815223409615
$ws.Range("E21").Value = 0.04514435695850949
$ws.Range("F21").Value = 2.416755148649116
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 1.669571627063135
$ws.Range("L21").Value = 0.3286610696631698

$ws.Range("B22").Value = 2.487011137446473
$ws.Range("C22").Value = 0.4451914768027336
$ws.Range("D22").Value = 0.08816136531699215
$ws.Range("E22").Value = 0.0450557071337192
$ws.Range("F22").Value = 2.47562109748074
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 1.694952996079579
$ws.Range("L22").Value = 0.3391936116701686

$ws.Range("B23").Value = 2.427113539550419
$ws.Range("C23").Value = 0.4290543171636045
$ws.Range("D23").Value = 0.0881726143326631
$ws.Range("E23").Value = 0.04510232933763714
$ws.Range("F23").Value = 2.444083688409393
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 1.681329631784195
$ws.Range("L23").Value = 0.3335569337926927

$ws.Range("B24").Value = 2.202347468076425
$ws.Range("C24").Value = 0.3680670378619766
$ws.Range("D24").Value = 0.08828169744065661
$ws.Range("E24").Value = 0.04529410462494976
$ws.Range("F24").Value = 2.327246930746554
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 1.631409990902156
$ws.Range("L24").Value = 0.3125425074752144

$ws.Range("B25").Value = 1.964299393445003
$ws.Range("C25").Value = 0.3025296680351062
$ws.Range("D25").Value = 0.08853838897888267
$ws.Range("E25").Value = 0.04553345119743391
$ws.Range("F25").Value = 2.206660695328992
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 1.58107224860801
$ws.Range("L25").Value = 0.29057716589476
